# Applies the cryptos-list price/volume refresh described by the commit.
# Each entry is (cell, new value, forceText) - forceText guards decimal-looking
# strings (e.g. '1.000', '0.7000') from Excel's automatic number coercion, which
# would otherwise silently drop the significant trailing zeros / type as Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = 'D2'; Value = '23.439.42'; ForceText = $false },
    @{ Ref = 'E2'; Value = '  -0.84%  '; ForceText = $false },
    @{ Ref = 'D3'; Value = '1.635.58'; ForceText = $false },
    @{ Ref = 'E3'; Value = '  -1.26%  '; ForceText = $false },
    @{ Ref = 'D4'; Value = '1.000'; ForceText = $true },
    @{ Ref = 'E4'; Value = '  +0.35%  '; ForceText = $false },
    @{ Ref = 'E5'; Value = '  +0.37%  '; ForceText = $false },
    @{ Ref = 'D6'; Value = '303.78'; ForceText = $true },
    @{ Ref = 'E6'; Value = '  -0.70%  '; ForceText = $false },
    @{ Ref = 'D7'; Value = '0.3784'; ForceText = $true },
    @{ Ref = 'E7'; Value = '  +0.35%  '; ForceText = $false },
    @{ Ref = 'D8'; Value = '51.87'; ForceText = $true },
    @{ Ref = 'E8'; Value = '  -0.42%  '; ForceText = $false },
    @{ Ref = 'D9'; Value = '0.3626'; ForceText = $true },
    @{ Ref = 'E9'; Value = '  -0.91%  '; ForceText = $false },
    @{ Ref = 'D10'; Value = '0.08180'; ForceText = $true },
    @{ Ref = 'E10'; Value = '  +0.20%  '; ForceText = $false },
    @{ Ref = 'D11'; Value = '1.228'; ForceText = $true },
    @{ Ref = 'E11'; Value = '  -3.22%  '; ForceText = $false },
    @{ Ref = 'D12'; Value = '1.001'; ForceText = $true },
    @{ Ref = 'D13'; Value = '22.43'; ForceText = $true },
    @{ Ref = 'E13'; Value = '  -3.23%  '; ForceText = $false },
    @{ Ref = 'D14'; Value = '6.475'; ForceText = $true },
    @{ Ref = 'E14'; Value = '  -3.37%  '; ForceText = $false },
    @{ Ref = 'D15'; Value = '7.386'; ForceText = $true },
    @{ Ref = 'E15'; Value = '  -0.13%  '; ForceText = $false },
    @{ Ref = 'E16'; Value = '  -3.05%  '; ForceText = $false },
    @{ Ref = 'D17'; Value = '1.629.25'; ForceText = $false },
    @{ Ref = 'E17'; Value = '  -1.42%  '; ForceText = $false },
    @{ Ref = 'D18'; Value = '94.93'; ForceText = $true },
    @{ Ref = 'E18'; Value = '  -0.54%  '; ForceText = $false },
    @{ Ref = 'D19'; Value = '0.06950'; ForceText = $true },
    @{ Ref = 'E19'; Value = '  +0.73%  '; ForceText = $false },
    @{ Ref = 'D20'; Value = '6.569'; ForceText = $true },
    @{ Ref = 'E20'; Value = '  -0.61%  '; ForceText = $false },
    @{ Ref = 'D21'; Value = '17.52'; ForceText = $true },
    @{ Ref = 'E21'; Value = '  -4.97%  '; ForceText = $false },
    @{ Ref = 'D22'; Value = '0.9999'; ForceText = $true },
    @{ Ref = 'E22'; Value = '  +0.33%  '; ForceText = $false },
    @{ Ref = 'D23'; Value = '12.55'; ForceText = $true },
    @{ Ref = 'E23'; Value = '  -2.91%  '; ForceText = $false },
    @{ Ref = 'D24'; Value = '23.445.53'; ForceText = $false },
    @{ Ref = 'E24'; Value = '  -0.82%  '; ForceText = $false },
    @{ Ref = 'D25'; Value = '2.513'; ForceText = $true },
    @{ Ref = 'E25'; Value = '  +4.41%  '; ForceText = $false },
    @{ Ref = 'D26'; Value = '3.047'; ForceText = $true },
    @{ Ref = 'E26'; Value = '  -3.40%  '; ForceText = $false },
    @{ Ref = 'D27'; Value = '21.13'; ForceText = $true },
    @{ Ref = 'E27'; Value = '  -1.24%  '; ForceText = $false },
    @{ Ref = 'D28'; Value = '150.65'; ForceText = $true },
    @{ Ref = 'E28'; Value = '  -0.28%  '; ForceText = $false },
    @{ Ref = 'D29'; Value = '5.282'; ForceText = $true },
    @{ Ref = 'E29'; Value = '  -0.71%  '; ForceText = $false },
    @{ Ref = 'D30'; Value = '133.21'; ForceText = $true },
    @{ Ref = 'E30'; Value = '  -2.77%  '; ForceText = $false },
    @{ Ref = 'D31'; Value = '1.811.84'; ForceText = $false },
    @{ Ref = 'E31'; Value = '  -1.37%  '; ForceText = $false },
    @{ Ref = 'D32'; Value = '6.631'; ForceText = $true },
    @{ Ref = 'E32'; Value = '  -4.10%  '; ForceText = $false },
    @{ Ref = 'D33'; Value = '2.164'; ForceText = $true },
    @{ Ref = 'E33'; Value = '  -6.85%  '; ForceText = $false },
    @{ Ref = 'D34'; Value = '1.046'; ForceText = $true },
    @{ Ref = 'E34'; Value = '  +7.44%  '; ForceText = $false },
    @{ Ref = 'D35'; Value = '11.27'; ForceText = $true },
    @{ Ref = 'E35'; Value = '  +1.95%  '; ForceText = $false },
    @{ Ref = 'D36'; Value = '0.02751'; ForceText = $true },
    @{ Ref = 'E36'; Value = '  -4.02%  '; ForceText = $false },
    @{ Ref = 'B37'; Value = 'Stellar'; ForceText = $false },
    @{ Ref = 'C37'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false },
    @{ Ref = 'D37'; Value = '0.08785'; ForceText = $true },
    @{ Ref = 'E37'; Value = '  -1.50%  '; ForceText = $false },
    @{ Ref = 'B38'; Value = 'Algorand'; ForceText = $false },
    @{ Ref = 'C38'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText = $false },
    @{ Ref = 'D38'; Value = '0.2494'; ForceText = $true },
    @{ Ref = 'E38'; Value = '  -3.22%  '; ForceText = $false },
    @{ Ref = 'D39'; Value = '0.07109'; ForceText = $true },
    @{ Ref = 'E39'; Value = '  -4.06%  '; ForceText = $false },
    @{ Ref = 'D40'; Value = '6.019'; ForceText = $true },
    @{ Ref = 'E40'; Value = '  -5.71%  '; ForceText = $false },
    @{ Ref = 'B41'; Value = 'TrustWalletToken'; ForceText = $false },
    @{ Ref = 'C41'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; ForceText = $false },
    @{ Ref = 'D41'; Value = '1.345'; ForceText = $true },
    @{ Ref = 'E41'; Value = '  -2.28%  '; ForceText = $false },
    @{ Ref = 'B42'; Value = 'TheSandbox'; ForceText = $false },
    @{ Ref = 'C42'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; ForceText = $false },
    @{ Ref = 'D42'; Value = '0.7000'; ForceText = $true },
    @{ Ref = 'E42'; Value = '  -2.26%  '; ForceText = $false },
    @{ Ref = 'D43'; Value = '15.82'; ForceText = $true },
    @{ Ref = 'E43'; Value = '  -3.93%  '; ForceText = $false },
    @{ Ref = 'D44'; Value = '12.14'; ForceText = $true },
    @{ Ref = 'E44'; Value = '  -3.70%  '; ForceText = $false },
    @{ Ref = 'D45'; Value = '0.6494'; ForceText = $true },
    @{ Ref = 'E45'; Value = '  -2.05%  '; ForceText = $false },
    @{ Ref = 'D46'; Value = '0.9998'; ForceText = $true },
    @{ Ref = 'E46'; Value = '  +0.36%  '; ForceText = $false },
    @{ Ref = 'D47'; Value = '2.272'; ForceText = $true },
    @{ Ref = 'E47'; Value = '  -3.87%  '; ForceText = $false },
    @{ Ref = 'D48'; Value = '3.968'; ForceText = $true },
    @{ Ref = 'E48'; Value = '  -1.44%  '; ForceText = $false },
    @{ Ref = 'D49'; Value = '0.07983'; ForceText = $true },
    @{ Ref = 'E49'; Value = '  -0.79%  '; ForceText = $false },
    @{ Ref = 'D50'; Value = '126.59'; ForceText = $true },
    @{ Ref = 'E50'; Value = '  -2.53%  '; ForceText = $false },
    @{ Ref = 'D51'; Value = '1.187'; ForceText = $true },
    @{ Ref = 'E51'; Value = '  -3.14%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Ref)
    if ($u.ForceText) {
        # Pre-format as Text so a purely-numeric-looking string (e.g. "303.78")
        # is stored verbatim instead of being parsed into a Number.
        $rng.NumberFormat = '@'
    }
    $rng.Value = $u.Value
    if ($u.ForceText) {
        # Drop back to the workbook default style so we do not leave a stray
        # explicit cell format behind (the source cells carry no style at all).
        $rng.Style = 'Normal'
    }
}
